$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 145 entirely (post removed), shifting subsequent rows up.
$ws.Rows.Item(145).Delete()
